$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on numeric-looking Price cells so they are stored as text
# (matching the original inline-string representation) instead of being auto-converted to numbers.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

# Apply the updated cell values from the crypto price refresh
$ws.Range('D2').Value = '24.579.36'
$ws.Range('E2').Value = '  +3.81%  '
$ws.Range('D3').Value = '1.694.12'
$ws.Range('E3').Value = '  +2.20%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '316.88'
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '0.3937'
$ws.Range('E7').Value = '  +1.45%  '
$ws.Range('D8').Value = '0.4018'
$ws.Range('E8').Value = '  +1.98%  '
$ws.Range('D9').Value = '1.532'
$ws.Range('E9').Value = '  +7.01%  '
$ws.Range('D10').Value = '54.07'
$ws.Range('E10').Value = '  +10.23%  '
$ws.Range('D11').Value = '1.001'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').Value = '0.08757'
$ws.Range('E12').Value = '  +1.36%  '
$ws.Range('D13').Value = '7.221'
$ws.Range('E13').Value = '  +8.12%  '
$ws.Range('E14').Value = '  +2.82%  '
$ws.Range('D15').Value = '0.00001320'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').Value = '7.606'
$ws.Range('E16').Value = '  +4.91%  '
$ws.Range('D17').Value = '1.699.05'
$ws.Range('E17').Value = '  +2.33%  '
$ws.Range('D18').Value = '100.88'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('D19').Value = '0.07035'
$ws.Range('E19').Value = '  +3.28%  '
$ws.Range('D20').Value = '19.62'
$ws.Range('E20').Value = '  +3.10%  '
$ws.Range('D21').Value = '6.851'
$ws.Range('E21').Value = '  +2.81%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = '14.04'
$ws.Range('E23').Value = '  +1.04%  '
$ws.Range('D24').Value = '24.578.06'
$ws.Range('E24').Value = '  +3.89%  '
$ws.Range('D25').Value = '3.029'
$ws.Range('E25').Value = '  +8.30%  '
$ws.Range('D26').Value = '2.315'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = '22.36'
$ws.Range('E27').Value = '  +2.94%  '
$ws.Range('D28').Value = '159.28'
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('D29').Value = '5.204'
$ws.Range('E29').Value = '  +1.10%  '
$ws.Range('D30').Value = '133.85'
$ws.Range('E30').Value = '  +3.19%  '
$ws.Range('D31').Value = '7.534'
$ws.Range('E31').Value = '  +15.67%  '
$ws.Range('D32').Value = '1.885.66'
$ws.Range('E32').Value = '  +2.38%  '
$ws.Range('D33').Value = '1.099'
$ws.Range('E33').Value = '  -2.83%  '
$ws.Range('D34').Value = '7.318'
$ws.Range('E34').Value = '  +12.41%  '
$ws.Range('D35').Value = '0.08526'
$ws.Range('E35').Value = '  -0.45%  '
$ws.Range('D36').Value = '11.40'
$ws.Range('E36').Value = '  +9.95%  '
$ws.Range('D37').Value = '1.975'
$ws.Range('E37').Value = '  +0.59%  '
$ws.Range('D38').Value = '0.2725'
$ws.Range('E38').Value = '  +3.01%  '
$ws.Range('D39').Value = '14.56'
$ws.Range('E39').Value = '  +0.74%  '
$ws.Range('D40').Value = '0.02749'
$ws.Range('E40').Value = '  +9.07%  '
$ws.Range('D41').Value = '0.09038'
$ws.Range('E41').Value = '  +2.95%  '
$ws.Range('D42').Value = '1.467'
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('D43').Value = '0.7690'
$ws.Range('E43').Value = '  +1.92%  '
$ws.Range('D44').Value = '0.7196'
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('D45').Value = '15.34'
$ws.Range('E45').Value = '  +3.47%  '
$ws.Range('D46').Value = '2.510'
$ws.Range('E46').Value = '  +4.55%  '
$ws.Range('D47').Value = '4.215'
$ws.Range('E47').Value = '  +2.99%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').Value = '1.001'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('B49').Value = 'Flow'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D49').Value = '1.351'
$ws.Range('E49').Value = '  +14.02%  '
$ws.Range('D50').Value = '141.30'
$ws.Range('E50').Value = '  +2.47%  '
$ws.Range('D51').Value = '0.08023'
$ws.Range('E51').Value = '  +2.99%  '
